$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.725.33'
$ws.Range('E2').Value = '  +7.78%  '
$ws.Range('D3').Value = '3.634.03'
$ws.Range('E3').Value = '  +7.49%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '593.08'
$ws.Range('E5').Value = '  +5.53%  '
$ws.Range('D6').Value = '192.30'
$ws.Range('E6').Value = '  +9.83%  '
$ws.Range('E7').Value = '  +3.08%  '
$ws.Range('D8').Value = '3.608.33'
$ws.Range('E8').Value = '  +6.87%  '
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('E10').Value = '  +5.24%  '
$ws.Range('D11').Value = '0.663'
$ws.Range('E11').Value = '  +4.73%  '
$ws.Range('D12').Value = '57.81'
$ws.Range('E12').Value = '  +8.31%  '
$ws.Range('D13').Value = '0.0000296'
$ws.Range('E13').Value = '  +6.97%  '
$ws.Range('E14').Value = '  +5.83%  '
$ws.Range('D15').Value = '4.214.36'
$ws.Range('E15').Value = '  +7.55%  '
$ws.Range('D16').Value = '3.628.95'
$ws.Range('E16').Value = '  +7.51%  '
$ws.Range('D17').Value = '19.36'
$ws.Range('E17').Value = '  +6.41%  '
$ws.Range('D18').Value = '70.501.31'
$ws.Range('E18').Value = '  +7.40%  '
$ws.Range('D19').Value = '12.60'
$ws.Range('E19').Value = '  +6.43%  '
$ws.Range('E20').Value = '  +1.09%  '
$ws.Range('E21').Value = '  +5.33%  '
$ws.Range('D22').Value = '495.69'
$ws.Range('E22').Value = '  +5.74%  '
$ws.Range('E23').Value = '  +13.09%  '
$ws.Range('D24').Value = '16.66'
$ws.Range('E24').Value = '  +16.48%  '
$ws.Range('D25').Value = '4.44'
$ws.Range('E25').Value = '  +8.71%  '
$ws.Range('D26').Value = '90.66'
$ws.Range('E26').Value = '  +0.57%  '
$ws.Range('D27').Value = '3.10'
$ws.Range('E27').Value = '  +6.62%  '
$ws.Range('D28').Value = '11.20'
$ws.Range('E28').Value = '  +5.77%  '
$ws.Range('D29').Value = '9.38'
$ws.Range('E29').Value = '  +7.66%  '
$ws.Range('D30').Value = '32.41'
$ws.Range('E30').Value = '  +4.10%  '
$ws.Range('D31').Value = '7.60'
$ws.Range('E31').Value = '  +15.54%  '
$ws.Range('E32').Value = '  +6.96%  '
$ws.Range('D33').Value = '614.39'
$ws.Range('E33').Value = '  +6.54%  '
$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D34').Value = '65.43'
$ws.Range('E34').Value = '  +5.50%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value = '0.116'
$ws.Range('E35').Value = '  +7.61%  '
$ws.Range('E36').Value = '  +12.77%  '
$ws.Range('E37').Value = '  +5.16%  '
$ws.Range('B38').Value = 'Dai'
$ws.Range('C38').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  +0.06%  '
$ws.Range('B39').Value = 'InjectiveProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D39').Value = '38.00'
$ws.Range('E39').Value = '  +6.09%  '
$ws.Range('B40').Value = 'TheGraph'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D40').Value = '0.402'
$ws.Range('E40').Value = '  +6.84%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').Value = '3.67'
$ws.Range('E41').Value = '  +2.34%  '
$ws.Range('D42').Value = '3.360.12'
$ws.Range('E42').Value = '  +8.37%  '
$ws.Range('D43').Value = '3.07'
$ws.Range('E43').Value = '  +8.70%  '
$ws.Range('E44').Value = '  +7.00%  '
$ws.Range('D45').Value = '2.68'
$ws.Range('E45').Value = '  +10.03%  '
$ws.Range('D46').Value = '3.35'
$ws.Range('E46').Value = '  +6.04%  '
$ws.Range('E47').Value = '  +2.91%  '
$ws.Range('E48').Value = '  +13.49%  '
$ws.Range('D49').Value = '9.09'
$ws.Range('E49').Value = '  +7.22%  '
$ws.Range('D50').Value = '3.29'
$ws.Range('E50').Value = '  +4.73%  '
$ws.Range('D51').Value = '0.999'
$ws.Range('E51').Value = '  -0.18%  '
